# "Add files via upload" — re-upload of a cleaned-up score sheet.
# The second worksheet (用戶行為分析及行銷決策(行銷三合)) contained nine
# duplicate / garbled student rows (typo'd IDs like "D114558t", a stray
# leading space on " 洪儀", a trailing "」" on "邱芊瑜」", etc.) that were
# removed from the source data. Delete those rows outright so everything
# below shifts up, exactly like the re-uploaded workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Row numbers refer to the *original* (pre-edit) layout, so delete from
# the bottom up to avoid invalidating not-yet-processed row numbers.
$rowsToDelete = @(69, 66, 64, 61, 54, 37, 35, 27, 7)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
